$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 20:35"

# --- Update Estados Unidos (row 4) stats ---
$ws.Range("B4").Value = 1826909
$ws.Range("C4").Value = 10089
$ws.Range("E4").Value = 1182436
$ws.Range("G4").Value = 329
$ws.Range("H4").Value = 105886

# --- Egipto had a big jump in cases and now ranks right after Irlanda,
#     ahead of Polonia and Ucrania. Re-sort that block of rows (38-40):
#     row 38 becomes Egipto (with fresh numbers), Polonia and Ucrania
#     both shift down one row, keeping their existing figures. ---
$ws.Range("A38").Value = "Egipto"
$ws.Range("B38").Value = 24985
$ws.Range("C38").Value = 1536
$ws.Range("D38").Value = 6810
$ws.Range("E38").Value = 17216
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 46
$ws.Range("H38").Value = 959

$ws.Range("A39").Value = "Polonia"
$ws.Range("B39").Value = 23786
$ws.Range("C39").Value = 215
$ws.Range("D39").Value = 11271
$ws.Range("E39").Value = 11451
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 1064

$ws.Range("A40").Value = "Ucrania"
$ws.Range("B40").Value = 23672
$ws.Range("C40").Value = 468
$ws.Range("D40").Value = 9538
$ws.Range("E40").Value = 13426
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 12
$ws.Range("H40").Value = 708
